$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 526.45
$ws.Range("C2").Value = 509.05
$ws.Range("D2").Value = 523
$ws.Range("E2").Value = 524.25
$ws.Range("F2").Value = 17
$ws.Range("G2").Value = 516.3
$ws.Range("B3").Value = 3031.7
$ws.Range("C3").Value = 2974.05
$ws.Range("D3").Value = 2999.95
$ws.Range("E3").Value = 3002
$ws.Range("F3").Value = 7
$ws.Range("G3").Value = 2979
$ws.Range("B4").Value = 510.9
$ws.Range("C4").Value = 498.2
$ws.Range("D4").Value = 505.55
$ws.Range("E4").Value = 505.75
$ws.Range("F4").Value = 18
$ws.Range("G4").Value = 499.15
$ws.Range("B5").Value = 1833.3
$ws.Range("C5").Value = 1800
$ws.Range("D5").Value = 1821.5
$ws.Range("E5").Value = 1819.9
$ws.Range("F5").Value = 9
$ws.Range("G5").Value = 1807
$ws.Range("B6").Value = 6938.65
$ws.Range("C6").Value = 6820.7
$ws.Range("D6").Value = 6910.05
$ws.Range("E6").Value = 6899.55
$ws.Range("F6").Value = 7
$ws.Range("G6").Value = 6844.25
$ws.Range("B7").Value = 192.68
$ws.Range("C7").Value = 188.68
$ws.Range("D7").Value = 191.96
$ws.Range("E7").Value = 191.96
$ws.Range("F7").Value = 98
$ws.Range("G7").Value = 188.9
$ws.Range("B8").Value = 248.75
$ws.Range("C8").Value = 237.85
$ws.Range("D8").Value = 247.5
$ws.Range("E8").Value = 247.97
$ws.Range("F8").Value = 109
$ws.Range("G8").Value = 238.2
$ws.Range("B9").Value = 493.75
$ws.Range("C9").Value = 483.35
$ws.Range("D9").Value = 491.05
$ws.Range("E9").Value = 492.2
$ws.Range("F9").Value = 45
$ws.Range("G9").Value = 484.7
$ws.Range("B10").Value = 880.35
$ws.Range("C10").Value = 848.7
$ws.Range("D10").Value = 874
$ws.Range("E10").Value = 875.15
$ws.Range("F10").Value = 24
$ws.Range("G10").Value = 851.25
$ws.Range("B11").Value = 4799.9
$ws.Range("C11").Value = 4670.6
$ws.Range("D11").Value = 4761.75
$ws.Range("E11").Value = 4765.65
$ws.Range("F11").Value = 7
$ws.Range("G11").Value = 4691.6
$ws.Range("B12").Value = 196.25
$ws.Range("C12").Value = 190.59
$ws.Range("D12").Value = 195.05
$ws.Range("E12").Value = 195.3
$ws.Range("F12").Value = 73
$ws.Range("G12").Value = 190.86
$ws.Range("B13").Value = 1870.5
$ws.Range("C13").Value = 1843
$ws.Range("D13").Value = 1856.55
$ws.Range("E13").Value = 1857.7
$ws.Range("F13").Value = 23
$ws.Range("G13").Value = 1862.95
$ws.Range("B14").Value = 1691
$ws.Range("C14").Value = 1654.1
$ws.Range("D14").Value = 1684.8
$ws.Range("E14").Value = 1681.85
$ws.Range("F14").Value = 124
$ws.Range("G14").Value = 1655.7
$ws.Range("B15").Value = 756.45
$ws.Range("C15").Value = 725.5
$ws.Range("D15").Value = 753.2
$ws.Range("E15").Value = 753.5
$ws.Range("F15").Value = 54
$ws.Range("G15").Value = 725.6
$ws.Range("B16").Value = 1268.9
$ws.Range("C16").Value = 1225.75
$ws.Range("D16").Value = 1267
$ws.Range("E16").Value = 1264.5
$ws.Range("F16").Value = 110
$ws.Range("G16").Value = 1228.55
$ws.Range("B17").Value = 1354.95
$ws.Range("C17").Value = 1331.5
$ws.Range("D17").Value = 1347.95
$ws.Range("E17").Value = 1347.25
$ws.Range("F17").Value = 36
$ws.Range("G17").Value = 1333.85
$ws.Range("B18").Value = 1936.85
$ws.Range("C18").Value = 1869.25
$ws.Range("D18").Value = 1885
$ws.Range("E18").Value = 1879.6
$ws.Range("F18").Value = 105
$ws.Range("G18").Value = 1930
$ws.Range("B19").Value = 973.2
$ws.Range("C19").Value = 923.7
$ws.Range("D19").Value = 965
$ws.Range("E19").Value = 964.5
$ws.Range("F19").Value = 26
$ws.Range("G19").Value = 924
$ws.Range("B20").Value = 619.9
$ws.Range("C20").Value = 602.25
$ws.Range("D20").Value = 619.6
$ws.Range("E20").Value = 618.6
$ws.Range("F20").Value = 9
$ws.Range("G20").Value = 602.55
$ws.Range("B21").Value = 2988.85
$ws.Range("C21").Value = 2928
$ws.Range("D21").Value = 2965
$ws.Range("E21").Value = 2964.25
$ws.Range("F21").Value = 34
$ws.Range("G21").Value = 2940.75
$ws.Range("B22").Value = 293.4
$ws.Range("C22").Value = 280
$ws.Range("D22").Value = 291.1
$ws.Range("E22").Value = 291
$ws.Range("F22").Value = 34
$ws.Range("G22").Value = 280.5
$ws.Range("B23").Value = 425.95
$ws.Range("C23").Value = 414.85
$ws.Range("D23").Value = 425.1
$ws.Range("E23").Value = 424.95
$ws.Range("F23").Value = 88
$ws.Range("G23").Value = 414.95
$ws.Range("B24").Value = 2736.85
$ws.Range("C24").Value = 2682.65
$ws.Range("D24").Value = 2717
$ws.Range("E24").Value = 2718.6
$ws.Range("F24").Value = 60
$ws.Range("G24").Value = 2684.1
$ws.Range("B25").Value = 822.6
$ws.Range("C25").Value = 802.9
$ws.Range("D25").Value = 819.95
$ws.Range("E25").Value = 820.4
$ws.Range("F25").Value = 117
$ws.Range("G25").Value = 804.95
$ws.Range("B26").Value = 763.95
$ws.Range("C26").Value = 757
$ws.Range("D26").Value = 759.15
$ws.Range("E26").Value = 759.95
$ws.Range("F26").Value = 4
$ws.Range("G26").Value = 757.3
$ws.Range("B27").Value = 1093.45
$ws.Range("C27").Value = 1040
$ws.Range("D27").Value = 1083.95
$ws.Range("E27").Value = 1088.1
$ws.Range("F27").Value = 21
$ws.Range("G27").Value = 1047.95
$ws.Range("B28").Value = 917.7
$ws.Range("C28").Value = 889.2
$ws.Range("D28").Value = 911
$ws.Range("E28").Value = 910.15
$ws.Range("F28").Value = 115
$ws.Range("G28").Value = 890.5
$ws.Range("B29").Value = 456.5
$ws.Range("C29").Value = 443.15
$ws.Range("D29").Value = 453.55
$ws.Range("E29").Value = 453.55
$ws.Range("F29").Value = 110
$ws.Range("G29").Value = 443.7
$ws.Range("B30").Value = 156.2
$ws.Range("C30").Value = 151.1
$ws.Range("D30").Value = 155.5
$ws.Range("E30").Value = 155.39
$ws.Range("F30").Value = 347
$ws.Range("G30").Value = 151.16
$ws.Range("B31").Value = 11139.95
$ws.Range("C31").Value = 10901
$ws.Range("D31").Value = 11084
$ws.Range("E31").Value = 11069.3
$ws.Range("F31").Value = 3
$ws.Range("G31").Value = 10915.35
